$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2025-04-28 06:52:57"
$ws.Range("B5").Value = 202
